# The captured XML diff for this revision touches only the root-element
# namespace-prefix declarations (xmlns:m / xmlns:ns17->ns19 reordering) in
# document.xml, endnotes.xml, footer1.xml, footer2.xml, footnotes.xml,
# header1.xml, numbering.xml, styles.xml, theme1.xml and customXml/item1.xml.
# Every one of those hunks maps to exactly the same single opening-tag line;
# there is no change anywhere inside the body of any part (the hunk right
# after the document.xml root tag doesn't resume until line 1421, i.e. past
# the end of that 1419-line part). The set of namespace URIs declared is
# identical before and after - only the prefix spelling/ordering used by the
# serializer differs, which is a cosmetic artifact of whatever tool produced
# the "after" package, not a content edit. There is nothing in the Word
# object model to change here, so we simply touch the active document
# without altering any content.
$d = $word.ActiveDocument
$d.Content | Out-Null
